# Initial deployment of Azure VMs from Excel
# Adds a second VM row (row 3) mirroring row 2's layout/format, wires up
# its "Password" hyperlink, and updates the view selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 values (column order matches the header / row 2 layout) ---
$ws.Range("A3").Value = "83140706-7c33-427a-a373-27883c159e91"
$ws.Range("B3").Value = "jf-rg-001"
$ws.Range("C3").Value = "jf-vm-001"
$ws.Range("D3").Value = "UK South"
$ws.Range("F3").Value = "MicrosoftWindowsServer"
$ws.Range("G3").Value = "WindowsServer"
$ws.Range("H3").Value = "2012-R2-Datacenter"
$ws.Range("I3").Value = "Latest"
$ws.Range("J3").Value = "Standard_DS1_v2"
$ws.Range("K3").Value = "AdminUser"
$ws.Range("L3").Value = "P@ssword1234"
$ws.Range("M3").Value = "No"
$ws.Range("N3").Value = "jf-vnet-001"
$ws.Range("O3").Value = "default"
$ws.Range("P3").Value = "jf-vm-001754"
$ws.Range("W3").Value = "Premium_LRS"
$ws.Range("X3").Value = 127

# Row 3 should carry the same row height / formatting as row 2, so copy
# row 2's formats down before layering the hyperlink on L3.
$ws.Range("A2:AF2").Copy()
$ws.Range("A3:AF3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Hyperlink for the Password cell (mirrors L2's mailto hyperlink) ---
$ws.Hyperlinks.Add($ws.Range("L3"), "mailto:P@ssword1234")

# Re-apply the Hyperlink style/format to L3 since adding the hyperlink
# can touch formatting - keep it consistent with L2.
$ws.Range("L2").Copy()
$ws.Range("L3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- View state: selection moved to AH3, scrolled right to column U ---
$ws.Range("AH3").Select()
